# Generate Report for Handback
#
# For both the zh-cn and de-de target-language sheets, the two real
# (non ".localization-config") rows have been handed back:
#   - Status switches from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - The "Latest Target File" / "Latest Handback File" columns (E/F),
#     previously blank, are now filled in with hyperlinks mirroring the
#     "Source File Name" / "Latest Handoff File" columns (A/C)
#   - "Latest Handback DateTime" (G) moves from the sentinel
#     0001-01-01 00:00:00 to the actual handback timestamp

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# The "Overview" sheet's Status columns (B/C) mirror the same status text
# as the per-language sheets for the two real file rows.
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(2, 2).Value = $newStatus
$overview.Cells.Item(2, 3).Value = $newStatus
$overview.Cells.Item(3, 2).Value = $newStatus
$overview.Cells.Item(3, 3).Value = $newStatus

$sheets = @(
    @{
        Name = "zh-cn"
        HandbackTime = "2016-03-11 03:10:06"
        Rows = @(
            @{
                Row = 2
                MdTarget  = "https://github.com/OpenLocalizationTest/oltest/blob/598fbc58937aabf830798eeeba8f47d898b71838/e2e/0f05aa6a-6e67-49f4-8f9f-871630e02cef.md"
                XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e97d7a65dbe358395094cc2b9e2caa9d36b32ee7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0f05aa6a-6e67-49f4-8f9f-871630e02cef.b0f114c12b42e23bfbd56955c2abe60c53fa2f7f.zh-cn.xlf"
            },
            @{
                Row = 3
                MdTarget  = "https://github.com/OpenLocalizationTest/oltest/blob/598fbc58937aabf830798eeeba8f47d898b71838/e2e/eb5fca58-5717-4c65-9b3d-2ba88abb2acc.md"
                XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e97d7a65dbe358395094cc2b9e2caa9d36b32ee7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eb5fca58-5717-4c65-9b3d-2ba88abb2acc.513060a60ff2ea96e1572eaf5dc9e7f676869012.zh-cn.xlf"
            }
        )
    },
    @{
        Name = "de-de"
        HandbackTime = "2016-03-11 03:10:30"
        Rows = @(
            @{
                Row = 2
                MdTarget  = "https://github.com/OpenLocalizationTest/oltest/blob/598fbc58937aabf830798eeeba8f47d898b71838/e2e/0f05aa6a-6e67-49f4-8f9f-871630e02cef.md"
                XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5afb8360233fd0b95a65ed724720731560e55ad/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0f05aa6a-6e67-49f4-8f9f-871630e02cef.b0f114c12b42e23bfbd56955c2abe60c53fa2f7f.de-de.xlf"
            },
            @{
                Row = 3
                MdTarget  = "https://github.com/OpenLocalizationTest/oltest/blob/598fbc58937aabf830798eeeba8f47d898b71838/e2e/eb5fca58-5717-4c65-9b3d-2ba88abb2acc.md"
                XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5afb8360233fd0b95a65ed724720731560e55ad/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eb5fca58-5717-4c65-9b3d-2ba88abb2acc.513060a60ff2ea96e1572eaf5dc9e7f676869012.de-de.xlf"
            }
        )
    }
)

foreach ($sheetInfo in $sheets) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    foreach ($rowInfo in $sheetInfo.Rows) {
        $r = $rowInfo.Row

        $sourceName = $ws.Cells.Item($r, 1).Text
        $handoffFileName = $ws.Cells.Item($r, 3).Text

        # Status -> handed back
        $ws.Cells.Item($r, 2).Value = $newStatus

        # E = Latest Target File (mirrors A: source file, as a hyperlink)
        $eCell = $ws.Cells.Item($r, 5)
        $ws.Hyperlinks.Add($eCell, $rowInfo.MdTarget, "", "", $sourceName)

        # F = Latest Handback File (mirrors C: handoff/target xlf, as a hyperlink)
        $fCell = $ws.Cells.Item($r, 6)
        $ws.Hyperlinks.Add($fCell, $rowInfo.XlfTarget, "", "", $handoffFileName)

        # G = Latest Handback DateTime
        $ws.Cells.Item($r, 7).Value = $sheetInfo.HandbackTime
    }
}
